$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (E1) onto the new
# header cell (F1), then set its text. This reproduces the same style
# index (bold, bordered, centered) used by the other header cells.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("F1").Value = "Modelo"

# Add the corresponding data cell for row 2 (no special style, like the
# other data cells).
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
